# Rajasthan Royals / Anuj Rawat batting log:
#  - rename the sheet to the player's name
#  - insert a new leading "matchNo" column and populate it

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts the existing
# teamName..result columns from A:L to B:M.
$ws.Columns.Item(1).Insert()

# Populate the new matchNo column (header + the single data row).
$ws.Range("A1").Value = "matchNo"
$ws.Range("A2").Value = "54th"

# Rename the sheet from the generic "Sheet1" to the player's name.
$ws.Name = "Anuj Rawat"
